# Revise error message display
#
# Appends one new data row (row 91) to the end of each of the 4 worksheets.
# Each new row duplicates the last existing row (row 90) on its sheet, but
# with the timestamp in column A updated to the new reading's time.

$wb = $excel.ActiveWorkbook

$newTime = [double]"45877.4634375"

# --- Sheet "MID_LFT_#1" (sheet1) ---
$ws = $wb.Worksheets.Item("MID_LFT_#1")
$ws.Cells.Item(91, 1).Value = $newTime
$ws.Cells.Item(91, 1).NumberFormat = $ws.Cells.Item(90, 1).NumberFormat
$ws.Cells.Item(91, 2).Value = "0x01,0x90"
$ws.Cells.Item(91, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(91, 4).Value = "0x01,0x1C"
$ws.Cells.Item(91, 5).Value = "0x07"
$ws.Cells.Item(91, 6).Value = 400
$ws.Cells.Item(91, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(91, 8).Value = 284
$ws.Cells.Item(91, 9).Value = 7

# --- Sheet "MID_LFT_#2" (sheet2) ---
$ws = $wb.Worksheets.Item("MID_LFT_#2")
$ws.Cells.Item(91, 1).Value = $newTime
$ws.Cells.Item(91, 1).NumberFormat = $ws.Cells.Item(90, 1).NumberFormat
$ws.Cells.Item(91, 2).Value = "0x01,0x7c"
$ws.Cells.Item(91, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(91, 4).Value = "0x01,0x28"
$ws.Cells.Item(91, 5).Value = "0x19"
$ws.Cells.Item(91, 6).Value = 380
$ws.Cells.Item(91, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item(91, 8).Value = 296
$ws.Cells.Item(91, 9).Value = 25

# --- Sheet "MID_PLT_#1" (sheet3) ---
$ws = $wb.Worksheets.Item("MID_PLT_#1")
$ws.Cells.Item(91, 1).Value = $newTime
$ws.Cells.Item(91, 1).NumberFormat = $ws.Cells.Item(90, 1).NumberFormat
$ws.Cells.Item(91, 2).Value = "0x00,0x6e"
$ws.Cells.Item(91, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(91, 4).Value = "0x00,0x5F"
$ws.Cells.Item(91, 5).Value = "0x15"
$ws.Cells.Item(91, 6).Value = 110
$ws.Cells.Item(91, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(91, 8).Value = 95
$ws.Cells.Item(91, 9).Value = 15

# --- Sheet "MID_PLT_#2" (sheet4) ---
$ws = $wb.Worksheets.Item("MID_PLT_#2")
$ws.Cells.Item(91, 1).Value = $newTime
$ws.Cells.Item(91, 1).NumberFormat = $ws.Cells.Item(90, 1).NumberFormat
$ws.Cells.Item(91, 2).Value = "0x00,0x82"
$ws.Cells.Item(91, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(91, 4).Value = "0x00,0x75"
$ws.Cells.Item(91, 5).Value = "0x9"
$ws.Cells.Item(91, 6).Value = 130
$ws.Cells.Item(91, 7).Value = [double]"5.68631262647113e+23"
$ws.Cells.Item(91, 8).Value = 117
$ws.Cells.Item(91, 9).Value = 9
